# Updated symbol list on Fri Jan  6 20:25:48 UTC 2023 with GitHub Actions
# Refreshes Price / Volume(1h) figures for the crypto list, and fixes the
# BOLO / CoinbaseStockToken row ordering (rows 47 and 48 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values are written with a leading apostrophe so Excel
# keeps them as text (matching the workbook's inlineStr/text cells)
# instead of auto-converting to a number/percentage, and the style is
# reset to "Normal" afterwards so no stray quote-prefix formatting is
# left behind on the cell.
$ws.Range("D2").Value = "'258.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.54%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.83%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.716"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.55%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06043"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.67%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'0.44%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8584"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.22%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9328"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.99%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1397"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.85%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05022"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'27.86%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07056"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03129"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.44%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09138"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.42%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001545"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.36%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006047"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.33%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005992"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.68%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.465"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.22%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.160"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.41%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'0.3094"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.31%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.28%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.115"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.51%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04232"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.24%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001213"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.77%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004046"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-5.79%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.25%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-21.51%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03861"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.27%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.31%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.003938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.02%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.01528"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'29.68%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002289"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.21%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.23%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.20%"
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.05442"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-9.26%"
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.1308"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-3.30%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.20%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("E50").Style = "Normal"
